$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Reuse the existing header style (same style index as H1) by copying formats
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New data values for columns I (I0) and J (IF), rows 2-6
$ws.Cells.Item(2, 9).Value = 7
$ws.Cells.Item(2, 10).Value = 8

$ws.Cells.Item(3, 9).Value = 6
$ws.Cells.Item(3, 10).Value = 7

$ws.Cells.Item(4, 9).Value = 7
$ws.Cells.Item(4, 10).Value = 7

$ws.Cells.Item(5, 9).Value = 2
$ws.Cells.Item(5, 10).Value = 3

$ws.Cells.Item(6, 9).Value = 7
$ws.Cells.Item(6, 10).Value = 7
